# Apply cyclic rotation of record data across rows 2-5 and a swap across
# rows 16-17 on the "Artfynd" sheet, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state per row, keyed by column letter. $null means the cell must
# end up empty (no value at all).
$targets = @{
    2 = @{
        A = 130826784; B = 57884; E = 100109; F = "Tretåig hackspett";
        G = "Picoides tridactylus"; H = "(Linnaeus, 1758)";
        M = "färska spår";
        P = "Brännan, Kälom, Offerdal, Jmt"; Q = 461233; R = 7039438; S = 10;
        Z = "11:37"; AB = "11:37"; AC = "Födosök barkfläk"
    }
    3 = @{
        A = 130825823; B = 57881; E = 100049; F = "Spillkråka";
        G = "Dryocopus martius"; H = "(Linnaeus, 1758)";
        M = "äldre spår";
        P = "Flinktorpet, Kälom, Offerdal, Jmt"; Q = 460947; R = 7039711; S = 10;
        Z = "10:38"; AB = "10:38"; AC = "Födosökshål på äldre döende gran."
    }
    4 = @{
        A = 130825852; B = 57884; E = 100109; F = "Tretåig hackspett";
        G = "Picoides tridactylus"; H = "(Linnaeus, 1758)";
        M = "färska spår";
        P = "Flinktorpet, Kälom, Offerdal, Jmt"; Q = 460952; R = 7039723; S = 15;
        Z = "10:42"; AB = "10:42"; AC = "Barkfläkta grövre och klenare granar."
    }
    5 = @{
        A = 130826010; B = 91808; E = 1202; F = "Ullticka";
        G = "Phellinidium ferrugineofuscum"; H = "(P.Karst.) Fiasson & Niemelä";
        M = $null;
        P = "Flinktorpet, Flinktorpet, Jmt"; Q = 460971; R = 7039688; S = 10;
        Z = "10:47"; AB = "10:47"; AC = $null
    }
    16 = @{
        A = 130826291; B = 91808; E = 1202; F = "Ullticka";
        G = "Phellinidium ferrugineofuscum"; H = "(P.Karst.) Fiasson & Niemelä";
        P = "Flinktorpet, Flinktorpet, Jmt"; Q = 461106; R = 7039672; S = 10;
        Z = "11:04"; AB = "11:04"; AC = $null
    }
    17 = @{
        A = 130825822; B = 89193; E = 510; F = "Doftskinn";
        G = "Cystostereum murrayi"; H = "(Berk. & M.A.Curtis.) Pouzar";
        P = "Flinktorpet, Flinktorpet, Jmt"; Q = 460947; R = 7039711; S = 10;
        Z = "10:36"; AB = "10:36"; AC = "På granlåga"
    }
}

foreach ($row in $targets.Keys) {
    $cols = $targets[$row]
    foreach ($col in $cols.Keys) {
        $value = $cols[$col]
        $addr = "$col$row"
        if ($null -eq $value) {
            $ws.Range($addr).Value = ""
        } else {
            $ws.Range($addr).Value = $value
        }
    }
}
